# cryptos.xlsx refresh — GitHub Actions scheduled price/volume update.
# For each data row (2-51) on the active sheet, write the latest "Price" (column D)
# and "Volume(1h)" (column E) values scraped for that coin.
#
# Column D values are written as literal text (not re-parsed as numbers) so that
# formatted figures such as "29.740.09", "10.30" or "0.00001037" survive exactly as
# scraped, matching how the source data feed renders them. NumberFormat is set to
# "@" (Text) right before the write and the cell style is restored to "Normal"
# immediately after, so the cell keeps its original (default) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Price = '29.740.09'; Volume = '  -0.05%  ' }
    @{ Row = 3; Price = '1.927.80'; Volume = '  -1.00%  ' }
    @{ Row = 4; Price = '0.9956'; Volume = '  -0.60%  ' }
    @{ Row = 5; Price = '334.86'; Volume = '  -1.95%  ' }
    @{ Row = 6; Price = '0.9966'; Volume = '  -0.46%  ' }
    @{ Row = 7; Price = '0.4673'; Volume = '  -2.36%  ' }
    @{ Row = 8; Price = '0.4168'; Volume = '  +0.97%  ' }
    @{ Row = 9; Price = '48.25'; Volume = '  +0.84%  ' }
    @{ Row = 10; Price = '0.08073'; Volume = '  -1.76%  ' }
    @{ Row = 11; Price = '1.027'; Volume = '  -0.78%  ' }
    @{ Row = 12; Price = '22.42'; Volume = '  -1.40%  ' }
    @{ Row = 13; Price = '1.934.32'; Volume = '  -1.28%  ' }
    @{ Row = 14; Price = '6.023'; Volume = '  -2.00%  ' }
    @{ Row = 15; Price = '7.195'; Volume = '  -2.41%  ' }
    @{ Row = 16; Price = '89.91'; Volume = '  -2.05%  ' }
    @{ Row = 17; Price = '0.9957'; Volume = '  -0.61%  ' }
    @{ Row = 18; Price = '0.00001037'; Volume = '  -1.88%  ' }
    @{ Row = 19; Price = '0.06597'; Volume = '  -1.09%  ' }
    @{ Row = 20; Price = '17.86'; Volume = '  -0.91%  ' }
    @{ Row = 21; Price = '0.9969'; Volume = '  -0.36%  ' }
    @{ Row = 22; Price = '29.673.84'; Volume = '  -0.17%  ' }
    @{ Row = 23; Price = '5.553'; Volume = '  -0.47%  ' }
    @{ Row = 24; Price = '11.53'; Volume = '  +2.60%  ' }
    @{ Row = 25; Price = '2.211'; Volume = '  -3.45%  ' }
    @{ Row = 26; Price = '2.146.67'; Volume = '  -1.79%  ' }
    @{ Row = 27; Price = '156.96'; Volume = '  -2.63%  ' }
    @{ Row = 28; Price = '19.97'; Volume = '  -1.14%  ' }
    @{ Row = 29; Price = '2.171'; Volume = '  +0.06%  ' }
    @{ Row = 30; Price = '5.685'; Volume = '  +0.69%  ' }
    @{ Row = 31; Price = '117.88'; Volume = '  -4.01%  ' }
    @{ Row = 32; Price = '1.044'; Volume = '  +3.78%  ' }
    @{ Row = 33; Price = '0.09474'; Volume = '  -1.80%  ' }
    @{ Row = 34; Price = '1.446'; Volume = '  -1.81%  ' }
    @{ Row = 35; Price = '5.459'; Volume = '  -0.69%  ' }
    @{ Row = 36; Price = '3.540'; Volume = '  -4.02%  ' }
    @{ Row = 37; Price = '0.06157'; Volume = '  -1.64%  ' }
    @{ Row = 38; Price = '0.02271'; Volume = '  -1.88%  ' }
    @{ Row = 39; Price = '8.472'; Volume = '  -0.36%  ' }
    @{ Row = 40; Price = '1.181'; Volume = '  -0.45%  ' }
    @{ Row = 41; Price = '0.5937'; Volume = '  -2.24%  ' }
    @{ Row = 42; Price = '0.9969'; Volume = '  -0.37%  ' }
    @{ Row = 43; Price = '10.30'; Volume = '  -3.75%  ' }
    @{ Row = 44; Price = '0.1843'; Volume = '  -2.63%  ' }
    @{ Row = 45; Price = '2.369'; Volume = '  -1.02%  ' }
    @{ Row = 46; Price = '1.259'; Volume = '  -1.15%  ' }
    @{ Row = 47; Price = '0.07563'; Volume = '  +1.97%  ' }
    @{ Row = 48; Price = '0.5611'; Volume = '  -1.76%  ' }
    @{ Row = 49; Price = '12.26'; Volume = '  -1.92%  ' }
    @{ Row = 50; Price = '1.946'; Volume = '  -2.09%  ' }
    @{ Row = 51; Price = '113.03'; Volume = '  +0.03%  ' }
)

foreach ($u in $updates) {
    $priceCell = $ws.Range("D$($u.Row)")
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $u.Price
    $priceCell.Style = "Normal"

    $ws.Range("E$($u.Row)").Value = $u.Volume
}
